# Generate Report for Handoff
# Regenerates the localization-status report: the Overview sheet and the two
# per-locale detail sheets (zh-cn, de-de) now reflect a new CI run against
# the callee/caller markdown fixtures instead of the old png/md fixtures,
# and a fourth source file (callerMd2.md) was added.

$wb = $excel.ActiveWorkbook

$HYPERLINK_UNDERLINE = $true
# Cornflower blue (FF6495ED) expressed as a VBA-style BGR integer so the
# Font.Color assignment below reproduces the workbook's existing custom
# "HyperLink" look instead of Excel's default theme hyperlink color.
$HYPERLINK_COLOR = 15570276

function Set-LinkCell($ws, $cellRef, $text, $url) {
    $ws.Range($cellRef).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
    $ws.Range($cellRef).Font.Underline = $HYPERLINK_UNDERLINE
    $ws.Range($cellRef).Font.Color = $HYPERLINK_COLOR
}

function Set-DateCell($ws, $cellRef, $text) {
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

Set-LinkCell $ov "A2" "calleeMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd1.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
Set-DateCell $ov "D2" "2016-03-22 04:59:18"

Set-LinkCell $ov "A3" "calleeMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd2.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
Set-DateCell $ov "D3" "2016-03-22 04:59:18"

Set-LinkCell $ov "A4" "callerMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd1.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
Set-DateCell $ov "D4" "2016-03-22 04:59:18"

Set-LinkCell $ov "A5" "callerMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd2.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
Set-DateCell $ov "D5" "2016-03-22 04:59:18"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

Set-LinkCell $zh "A2" "calleeMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd1.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
Set-LinkCell $zh "D2" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/307a54467126100e794936f4ad9fdf49dedb6de0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
Set-DateCell $zh "E2" "2016-03-22 04:59:14"
Set-DateCell $zh "H2" "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"
$zh.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

Set-LinkCell $zh "A3" "calleeMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd2.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
Set-LinkCell $zh "D3" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/307a54467126100e794936f4ad9fdf49dedb6de0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
Set-DateCell $zh "E3" "2016-03-22 04:59:14"
Set-DateCell $zh "H3" "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"
$zh.Range("K3").Value = "e2e\callerMd1.md"

Set-LinkCell $zh "A4" "callerMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd1.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
Set-LinkCell $zh "D4" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/307a54467126100e794936f4ad9fdf49dedb6de0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
Set-DateCell $zh "E4" "2016-03-22 04:59:14"
Set-DateCell $zh "H4" "0001-01-01 00:00:00"
$zh.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$zh.Range("J4").Value = "Include"

Set-LinkCell $zh "A5" "callerMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd2.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
Set-LinkCell $zh "D5" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/307a54467126100e794936f4ad9fdf49dedb6de0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
Set-DateCell $zh "E5" "2016-03-22 04:59:14"
Set-DateCell $zh "H5" "0001-01-01 00:00:00"
$zh.Range("I5").Value = "e2e\calleeMd1.md"
$zh.Range("J5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

Set-LinkCell $de "A2" "calleeMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd1.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
Set-LinkCell $de "D2" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86e5ef5a46704e792f4cec290bea8cae68eec63b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
Set-DateCell $de "E2" "2016-03-22 04:59:18"
Set-DateCell $de "H2" "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"
$de.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

Set-LinkCell $de "A3" "calleeMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/calleeMd2.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
Set-LinkCell $de "D3" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86e5ef5a46704e792f4cec290bea8cae68eec63b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
Set-DateCell $de "E3" "2016-03-22 04:59:18"
Set-DateCell $de "H3" "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"
$de.Range("K3").Value = "e2e\callerMd1.md"

Set-LinkCell $de "A4" "callerMd1.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd1.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
Set-LinkCell $de "D4" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86e5ef5a46704e792f4cec290bea8cae68eec63b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
Set-DateCell $de "E4" "2016-03-22 04:59:18"
Set-DateCell $de "H4" "0001-01-01 00:00:00"
$de.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$de.Range("J4").Value = "Include"

Set-LinkCell $de "A5" "callerMd2.md" "https://github.com/OpenLocalizationTest/oltest/blob/19fe6168dc2413cb9b03f7b310a5ec40c24c2b37/e2e/callerMd2.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
Set-LinkCell $de "D5" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86e5ef5a46704e792f4cec290bea8cae68eec63b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
Set-DateCell $de "E5" "2016-03-22 04:59:18"
Set-DateCell $de "H5" "0001-01-01 00:00:00"
$de.Range("I5").Value = "e2e\calleeMd1.md"
$de.Range("J5").Value = "Include"
